$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set column H (Industries) to 0 for rows 35 through 176
$ws.Range("H35:H176").Value = 0
